# Auto-generated edit script applying numeric updates per the commit diff.
# Each sheet is addressed by its fixed index (1-based) matching the workbook's tab order:
# 1=ALC 2=ARM 3=BSM 4=CRP 5=CUL 6=GSM 7=LTW 8=WVR
$wb = $excel.ActiveWorkbook

# --- Sheet 1 (ALC) ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(6, 8).Value2 = 1641.7059  # H6: 1453.8667 -> 1641.7059
$ws.Cells.Item(6, 10).Value2 = 2984.4443  # J6: 2965.5715 -> 2984.4443
$ws.Cells.Item(6, 12).Value2 = 8953.332900000001  # L6: 8896.7145 -> 8953.332900000001
$ws.Cells.Item(6, 14).Value2 = -9177.332900000001  # N6: -9120.7145 -> -9177.332900000001
$ws.Cells.Item(28, 8).Value2 = 465.41666  # H28: 433.85715 -> 465.41666
$ws.Cells.Item(28, 9).Value2 = 465.41666  # I28: 433.85715 -> 465.41666
$ws.Cells.Item(28, 11).Value2 = 465.41666  # K28: 433.85715 -> 465.41666
$ws.Cells.Item(28, 13).Value2 = 19.58334000000002  # M28: 51.14285000000001 -> 19.58334000000002
$ws.Cells.Item(38, 8).Value2 = 899  # H38: 899.46155 -> 899
$ws.Cells.Item(38, 9).Value2 = 199.6842  # I38: 199.78947 -> 199.6842
$ws.Cells.Item(38, 10).Value2 = 2797.1428  # J38: 2798.5715 -> 2797.1428
$ws.Cells.Item(38, 11).Value2 = 599.0526  # K38: 599.36841 -> 599.0526
$ws.Cells.Item(38, 12).Value2 = 8391.428400000001  # L38: 8395.7145 -> 8391.428400000001
$ws.Cells.Item(38, 13).Value2 = -227.0526  # M38: -227.36841 -> -227.0526
$ws.Cells.Item(38, 14).Value2 = -9135.428400000001  # N38: -9139.7145 -> -9135.428400000001
$ws.Cells.Item(98, 8).Value2 = 2087.875  # H98: 2267.9333 -> 2087.875
$ws.Cells.Item(98, 9).Value2 = 2071.3333  # I98: 2233.5454 -> 2071.3333
$ws.Cells.Item(98, 10).Value2 = 2137.5  # J98: 2362.5 -> 2137.5
$ws.Cells.Item(98, 11).Value2 = 2071.3333  # K98: 2233.5454 -> 2071.3333
$ws.Cells.Item(98, 12).Value2 = 2137.5  # L98: 2362.5 -> 2137.5
$ws.Cells.Item(98, 13).Value2 = -573.3332999999998  # M98: -735.5454 -> -573.3332999999998
$ws.Cells.Item(98, 14).Value2 = -5133.5  # N98: -5358.5 -> -5133.5
$ws.Cells.Item(122, 8).Value2 = 2087.875  # H122: 2267.9333 -> 2087.875
$ws.Cells.Item(122, 9).Value2 = 2071.3333  # I122: 2233.5454 -> 2071.3333
$ws.Cells.Item(122, 10).Value2 = 2137.5  # J122: 2362.5 -> 2137.5
$ws.Cells.Item(122, 11).Value2 = 6213.999899999999  # K122: 6700.6362 -> 6213.999899999999
$ws.Cells.Item(122, 12).Value2 = 6412.5  # L122: 7087.5 -> 6412.5
$ws.Cells.Item(122, 13).Value2 = -3763.999899999999  # M122: -4250.6362 -> -3763.999899999999
$ws.Cells.Item(122, 14).Value2 = -11312.5  # N122: -11987.5 -> -11312.5
$ws.Cells.Item(129, 8).Value2 = 945.2222  # H129: 930.3143 -> 945.2222
$ws.Cells.Item(129, 10).Value2 = 902.569  # J129: 891.10767 -> 902.569
$ws.Cells.Item(129, 12).Value2 = 2707.707  # L129: 2673.32301 -> 2707.707
$ws.Cells.Item(129, 14).Value2 = -12707.707  # N129: -12673.32301 -> -12707.707
$ws.Cells.Item(132, 8).Value2 = 4005997.5  # H132: 4451109.5 -> 4005997.5
$ws.Cells.Item(132, 9).Value2 = 4450589  # I132: 5135245 -> 4450589
$ws.Cells.Item(132, 10).Value2 = 4673.2  # J132: 4228.6665 -> 4673.2
$ws.Cells.Item(132, 11).Value2 = 13351767  # K132: 15405735 -> 13351767
$ws.Cells.Item(132, 12).Value2 = 14019.6  # L132: 12685.9995 -> 14019.6
$ws.Cells.Item(132, 13).Value2 = -13349237  # M132: -15403205 -> -13349237
$ws.Cells.Item(132, 14).Value2 = -19079.6  # N132: -17745.9995 -> -19079.6
$ws.Cells.Item(133, 8).Value2 = 20415.934  # H133: 22401.117 -> 20415.934
$ws.Cells.Item(133, 10).Value2 = 20415.934  # J133: 22401.117 -> 20415.934
$ws.Cells.Item(133, 12).Value2 = 20415.934  # L133: 22401.117 -> 20415.934
$ws.Cells.Item(133, 14).Value2 = -30535.934  # N133: -32521.117 -> -30535.934
$ws.Cells.Item(134, 8).Value2 = 25642.857  # H134: 26111.428 -> 25642.857
$ws.Cells.Item(134, 10).Value2 = 25642.857  # J134: 26111.428 -> 25642.857
$ws.Cells.Item(134, 12).Value2 = 25642.857  # L134: 26111.428 -> 25642.857
$ws.Cells.Item(134, 14).Value2 = -35782.857  # N134: -36251.428 -> -35782.857
$ws.Cells.Item(135, 8).Value2 = 1092.0588  # H135: 1348.625 -> 1092.0588
$ws.Cells.Item(135, 9).Value2 = 972.8125  # I135: 1282.9231 -> 972.8125
$ws.Cells.Item(135, 10).Value2 = 3000  # J135: 1633.3334 -> 3000
$ws.Cells.Item(135, 11).Value2 = 8755.3125  # K135: 11546.3079 -> 8755.3125
$ws.Cells.Item(135, 12).Value2 = 27000  # L135: 14700.0006 -> 27000
$ws.Cells.Item(135, 13).Value2 = -6220.3125  # M135: -9011.3079 -> -6220.3125
$ws.Cells.Item(135, 14).Value2 = -32070  # N135: -19770.0006 -> -32070
$ws.Cells.Item(137, 8).Value2 = 3128.8447  # H137: 3905.0466 -> 3128.8447
$ws.Cells.Item(137, 9).Value2 = 3139.976  # I137: 4647.72 -> 3139.976
$ws.Cells.Item(137, 10).Value2 = 3099.625  # J137: 2873.5557 -> 3099.625
$ws.Cells.Item(137, 11).Value2 = 9419.928  # K137: 13943.16 -> 9419.928
$ws.Cells.Item(137, 12).Value2 = 9298.875  # L137: 8620.667099999999 -> 9298.875
$ws.Cells.Item(137, 13).Value2 = -6869.928  # M137: -11393.16 -> -6869.928
$ws.Cells.Item(137, 14).Value2 = -14398.875  # N137: -13720.6671 -> -14398.875
$ws.Cells.Item(138, 8).Value2 = 1948.0674  # H138: 1939.6111 -> 1948.0674
$ws.Cells.Item(138, 9).Value2 = 1334.1794  # I138: 1335.9744 -> 1334.1794
$ws.Cells.Item(138, 10).Value2 = 2426.9  # J138: 2401.2156 -> 2426.9
$ws.Cells.Item(138, 11).Value2 = 4002.5382  # K138: 4007.9232 -> 4002.5382
$ws.Cells.Item(138, 12).Value2 = 7280.700000000001  # L138: 7203.6468 -> 7280.700000000001
$ws.Cells.Item(138, 13).Value2 = 1137.4618  # M138: 1132.0768 -> 1137.4618
$ws.Cells.Item(138, 14).Value2 = -17560.7  # N138: -17483.6468 -> -17560.7
$ws.Cells.Item(139, 8).Value2 = 29541.666  # H139: 29714.285 -> 29541.666
$ws.Cells.Item(139, 10).Value2 = 29541.666  # J139: 29714.285 -> 29541.666
$ws.Cells.Item(139, 12).Value2 = 29541.666  # L139: 29714.285 -> 29541.666
$ws.Cells.Item(139, 14).Value2 = -39821.666  # N139: -39994.285 -> -39821.666
$ws.Cells.Item(141, 8).Value2 = 803257.7  # H141: 714714.25 -> 803257.7
$ws.Cells.Item(141, 9).Value2 = 1994.7142  # I141: 2206.25 -> 1994.7142
$ws.Cells.Item(141, 10).Value2 = 6412098.5  # J141: 2139730.2 -> 6412098.5
$ws.Cells.Item(141, 11).Value2 = 5984.142599999999  # K141: 6618.75 -> 5984.142599999999
$ws.Cells.Item(141, 12).Value2 = 19236295.5  # L141: 6419190.600000001 -> 19236295.5
$ws.Cells.Item(141, 13).Value2 = -804.1425999999992  # M141: -1438.75 -> -804.1425999999992
$ws.Cells.Item(141, 14).Value2 = -19246655.5  # N141: -6429550.600000001 -> -19246655.5
# --- Sheet 2 (ARM) ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value2 = 7781.506  # H32: 8104.3765 -> 7781.506
$ws.Cells.Item(32, 9).Value2 = 6391.7437  # I32: 6687.4863 -> 6391.7437
$ws.Cells.Item(32, 11).Value2 = 6391.7437  # K32: 6687.4863 -> 6391.7437
$ws.Cells.Item(32, 13).Value2 = -6104.7437  # M32: -6400.4863 -> -6104.7437
$ws.Cells.Item(61, 8).Value2 = 5627.5835  # H61: 4807.7417 -> 5627.5835
$ws.Cells.Item(61, 9).Value2 = 5630.857  # I61: 6170.1577 -> 5630.857
$ws.Cells.Item(61, 10).Value2 = 5604.6665  # J61: 2650.5833 -> 5604.6665
$ws.Cells.Item(61, 11).Value2 = 5630.857  # K61: 6170.1577 -> 5630.857
$ws.Cells.Item(61, 12).Value2 = 5604.6665  # L61: 2650.5833 -> 5604.6665
$ws.Cells.Item(61, 13).Value2 = -5418.857  # M61: -5958.1577 -> -5418.857
$ws.Cells.Item(61, 14).Value2 = -6028.6665  # N61: -3074.5833 -> -6028.6665
$ws.Cells.Item(136, 8).Value2 = 5627.5835  # H136: 4807.7417 -> 5627.5835
$ws.Cells.Item(136, 9).Value2 = 5630.857  # I136: 6170.1577 -> 5630.857
$ws.Cells.Item(136, 10).Value2 = 5604.6665  # J136: 2650.5833 -> 5604.6665
$ws.Cells.Item(136, 11).Value2 = 16892.571  # K136: 18510.4731 -> 16892.571
$ws.Cells.Item(136, 12).Value2 = 16813.9995  # L136: 7951.749899999999 -> 16813.9995
$ws.Cells.Item(136, 13).Value2 = -14342.571  # M136: -15960.4731 -> -14342.571
$ws.Cells.Item(136, 14).Value2 = -21913.9995  # N136: -13051.7499 -> -21913.9995
# --- Sheet 3 (BSM) ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(99, 8).Value2 = 2760  # H99: 2908.3333 -> 2760
$ws.Cells.Item(99, 9).Value2 = 2052.5  # I99: 2275.7144 -> 2052.5
$ws.Cells.Item(99, 10).Value2 = 3274.5454  # J99: 3310.9092 -> 3274.5454
$ws.Cells.Item(99, 11).Value2 = 2052.5  # K99: 2275.7144 -> 2052.5
$ws.Cells.Item(99, 12).Value2 = 3274.5454  # L99: 3310.9092 -> 3274.5454
$ws.Cells.Item(99, 13).Value2 = -554.5  # M99: -777.7143999999998 -> -554.5
$ws.Cells.Item(99, 14).Value2 = -6270.5454  # N99: -6306.9092 -> -6270.5454
$ws.Cells.Item(134, 8).Value2 = 4114.4595  # H134: 4921.032 -> 4114.4595
$ws.Cells.Item(134, 9).Value2 = 3580.8708  # I134: 4236.96 -> 3580.8708
$ws.Cells.Item(134, 10).Value2 = 6871.3335  # J134: 7771.3335 -> 6871.3335
$ws.Cells.Item(134, 11).Value2 = 10742.6124  # K134: 12710.88 -> 10742.6124
$ws.Cells.Item(134, 12).Value2 = 20614.0005  # L134: 23314.0005 -> 20614.0005
$ws.Cells.Item(134, 13).Value2 = -8207.6124  # M134: -10175.88 -> -8207.6124
$ws.Cells.Item(134, 14).Value2 = -25684.0005  # N134: -28384.0005 -> -25684.0005
# --- Sheet 4 (CRP) ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(10, 8).Value2 = 5075.478  # H10: 20538.4 -> 5075.478
$ws.Cells.Item(10, 9).Value2 = 750  # I10: 607.1429000000001 -> 750
$ws.Cells.Item(10, 10).Value2 = 14962.286  # J10: 45905.453 -> 14962.286
$ws.Cells.Item(10, 11).Value2 = 750  # K10: 607.1429000000001 -> 750
$ws.Cells.Item(10, 12).Value2 = 14962.286  # L10: 45905.453 -> 14962.286
$ws.Cells.Item(10, 13).Value2 = -611  # M10: -468.1429000000001 -> -611
$ws.Cells.Item(10, 14).Value2 = -15240.286  # N10: -46183.453 -> -15240.286
$ws.Cells.Item(44, 8).Value2 = 0  # H44: 5000 -> 0
$ws.Cells.Item(44, 9).Value2 = 0  # I44: 5000 -> 0
$ws.Cells.Item(44, 11).Value2 = 0  # K44: 5000 -> 0
$ws.Cells.Item(44, 13).ClearContents()  # M44: -4558 -> (removed)
$ws.Cells.Item(52, 8).Value2 = 39450  # H52: 51300 -> 39450
$ws.Cells.Item(52, 10).Value2 = 39450  # J52: 51300 -> 39450
$ws.Cells.Item(52, 12).Value2 = 39450  # L52: 51300 -> 39450
$ws.Cells.Item(52, 14).Value2 = -40038  # N52: -51888 -> -40038
$ws.Cells.Item(99, 8).Value2 = 3200  # H99: 3750 -> 3200
$ws.Cells.Item(99, 10).Value2 = 4333.3335  # J99: 4500 -> 4333.3335
$ws.Cells.Item(99, 12).Value2 = 4333.3335  # L99: 4500 -> 4333.3335
$ws.Cells.Item(99, 14).Value2 = -7329.3335  # N99: -7496 -> -7329.3335
$ws.Cells.Item(126, 8).Value2 = 3200  # H126: 3750 -> 3200
$ws.Cells.Item(126, 10).Value2 = 4333.3335  # J126: 4500 -> 4333.3335
$ws.Cells.Item(126, 12).Value2 = 13000.0005  # L126: 13500 -> 13000.0005
$ws.Cells.Item(126, 14).Value2 = -17940.0005  # N126: -18440 -> -17940.0005
$ws.Cells.Item(132, 8).Value2 = 2066.8333  # H132: 2261.8386 -> 2066.8333
$ws.Cells.Item(132, 9).Value2 = 1496.6428  # I132: 1635.5217 -> 1496.6428
$ws.Cells.Item(132, 11).Value2 = 4489.928400000001  # K132: 4906.5651 -> 4489.928400000001
$ws.Cells.Item(132, 13).Value2 = -1959.928400000001  # M132: -2376.5651 -> -1959.928400000001
# --- Sheet 5 (CUL) ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(11, 8).Value2 = 34000  # H11: 25625 -> 34000
$ws.Cells.Item(11, 9).Value2 = 0  # I11: 500 -> 0
$ws.Cells.Item(11, 11).Value2 = 0  # K11: 1500 -> 0
$ws.Cells.Item(11, 13).ClearContents()  # M11: -1360 -> (removed)
$ws.Cells.Item(15, 8).Value2 = 1575  # H15: 1096.6666 -> 1575
$ws.Cells.Item(15, 9).Value2 = 150  # I15: 145 -> 150
$ws.Cells.Item(15, 11).Value2 = 450  # K15: 435 -> 450
$ws.Cells.Item(15, 13).Value2 = -310  # M15: -295 -> -310
$ws.Cells.Item(47, 8).Value2 = 1955.1111  # H47: 2281.6 -> 1955.1111
$ws.Cells.Item(47, 9).Value2 = 173  # I47: 69.666664 -> 173
$ws.Cells.Item(47, 10).Value2 = 3380.8  # J47: 3229.5715 -> 3380.8
$ws.Cells.Item(47, 11).Value2 = 519  # K47: 208.999992 -> 519
$ws.Cells.Item(47, 12).Value2 = 10142.4  # L47: 9688.7145 -> 10142.4
$ws.Cells.Item(47, 13).Value2 = -88  # M47: 222.000008 -> -88
$ws.Cells.Item(47, 14).Value2 = -11004.4  # N47: -10550.7145 -> -11004.4
$ws.Cells.Item(57, 8).Value2 = 4000  # H57: 0 -> 4000
$ws.Cells.Item(57, 10).Value2 = 4000  # J57: 0 -> 4000
$ws.Cells.Item(57, 12).Value2 = 12000  # L57: 0 -> 12000
$ws.Cells.Item(57, 14).Value2 = -13118  # N57: None -> -13118
$ws.Cells.Item(64, 8).Value2 = 2250.25  # H64: 11589.091 -> 2250.25
$ws.Cells.Item(64, 9).Value2 = 1500.5  # I64: 1745 -> 1500.5
$ws.Cells.Item(64, 10).Value2 = 3000  # J64: 13776.667 -> 3000
$ws.Cells.Item(64, 11).Value2 = 4501.5  # K64: 5235 -> 4501.5
$ws.Cells.Item(64, 12).Value2 = 9000  # L64: 41330.001 -> 9000
$ws.Cells.Item(64, 13).Value2 = -4231.5  # M64: -4965 -> -4231.5
$ws.Cells.Item(64, 14).Value2 = -9540  # N64: -41870.001 -> -9540
$ws.Cells.Item(67, 8).Value2 = 2250.25  # H67: 11589.091 -> 2250.25
$ws.Cells.Item(67, 9).Value2 = 1500.5  # I67: 1745 -> 1500.5
$ws.Cells.Item(67, 10).Value2 = 3000  # J67: 13776.667 -> 3000
$ws.Cells.Item(67, 11).Value2 = 4501.5  # K67: 5235 -> 4501.5
$ws.Cells.Item(67, 12).Value2 = 9000  # L67: 41330.001 -> 9000
$ws.Cells.Item(67, 13).Value2 = -3565.5  # M67: -4299 -> -3565.5
$ws.Cells.Item(67, 14).Value2 = -10872  # N67: -43202.001 -> -10872
# --- Sheet 6 (GSM) ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(122, 8).Value2 = 5201.5  # H122: 4921.3335 -> 5201.5
$ws.Cells.Item(122, 9).Value2 = 2962.6  # I122: 2635.3333 -> 2962.6
$ws.Cells.Item(122, 11).Value2 = 8887.799999999999  # K122: 7905.999899999999 -> 8887.799999999999
$ws.Cells.Item(122, 13).Value2 = -6437.799999999999  # M122: -5455.999899999999 -> -6437.799999999999
$ws.Cells.Item(124, 8).Value2 = 34745  # H124: 0 -> 34745
$ws.Cells.Item(124, 10).Value2 = 34745  # J124: 0 -> 34745
$ws.Cells.Item(124, 12).Value2 = 34745  # L124: 0 -> 34745
$ws.Cells.Item(124, 14).Value2 = -44565  # N124: None -> -44565
# --- Sheet 7 (LTW) ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(14, 8).Value2 = 10621.077  # H14: 372145.72 -> 10621.077
$ws.Cells.Item(14, 9).Value2 = 5502  # I14: 2500000 -> 5502
$ws.Cells.Item(14, 10).Value2 = 11551.818  # J14: 17503.334 -> 11551.818
$ws.Cells.Item(14, 11).Value2 = 5502  # K14: 2500000 -> 5502
$ws.Cells.Item(14, 12).Value2 = 11551.818  # L14: 17503.334 -> 11551.818
$ws.Cells.Item(14, 13).Value2 = -5330  # M14: -2499828 -> -5330
$ws.Cells.Item(14, 14).Value2 = -11895.818  # N14: -17847.334 -> -11895.818
$ws.Cells.Item(17, 8).Value2 = 33006  # H17: 32672.666 -> 33006
$ws.Cells.Item(17, 9).Value2 = 9000  # I17: 8000 -> 9000
$ws.Cells.Item(17, 11).Value2 = 9000  # K17: 8000 -> 9000
$ws.Cells.Item(17, 13).Value2 = -8830  # M17: -7830 -> -8830
$ws.Cells.Item(55, 8).Value2 = 768.8  # H55: 500.22223 -> 768.8
$ws.Cells.Item(55, 9).Value2 = 127.5  # I55: 188.33333 -> 127.5
$ws.Cells.Item(55, 10).Value2 = 1410.1  # J55: 1124 -> 1410.1
$ws.Cells.Item(55, 11).Value2 = 127.5  # K55: 188.33333 -> 127.5
$ws.Cells.Item(55, 12).Value2 = 1410.1  # L55: 1124 -> 1410.1
$ws.Cells.Item(55, 13).Value2 = 45.5  # M55: -15.33332999999999 -> 45.5
$ws.Cells.Item(55, 14).Value2 = -1756.1  # N55: -1470 -> -1756.1
$ws.Cells.Item(108, 8).Value2 = 31500  # H108: 32500 -> 31500
$ws.Cells.Item(108, 10).Value2 = 31500  # J108: 32500 -> 31500
$ws.Cells.Item(108, 12).Value2 = 31500  # L108: 32500 -> 31500
$ws.Cells.Item(108, 14).Value2 = -39180  # N108: -40180 -> -39180
$ws.Cells.Item(132, 8).Value2 = 3592.8096  # H132: 3682.15 -> 3592.8096
$ws.Cells.Item(132, 9).Value2 = 2454.4  # I132: 2484.3 -> 2454.4
$ws.Cells.Item(132, 10).Value2 = 4627.727  # J132: 4880 -> 4627.727
$ws.Cells.Item(132, 11).Value2 = 7363.200000000001  # K132: 7452.900000000001 -> 7363.200000000001
$ws.Cells.Item(132, 12).Value2 = 13883.181  # L132: 14640 -> 13883.181
$ws.Cells.Item(132, 13).Value2 = -4833.200000000001  # M132: -4922.900000000001 -> -4833.200000000001
$ws.Cells.Item(132, 14).Value2 = -18943.181  # N132: -19700 -> -18943.181
# --- Sheet 8 (WVR) ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(13, 8).Value2 = 20238.666  # H13: 56704 -> 20238.666
$ws.Cells.Item(13, 10).Value2 = 22756  # J13: 85006 -> 22756
$ws.Cells.Item(13, 12).Value2 = 22756  # L13: 85006 -> 22756
$ws.Cells.Item(13, 14).Value2 = -23036  # N13: -85286 -> -23036
$ws.Cells.Item(126, 8).Value2 = 2943481.5  # H126: 2704747 -> 2943481.5
$ws.Cells.Item(126, 9).Value2 = 1494.4615  # I126: 1358.8387 -> 1494.4615
$ws.Cells.Item(126, 10).Value2 = 12504939  # J126: 16672252 -> 12504939
$ws.Cells.Item(126, 11).Value2 = 4483.3845  # K126: 4076.5161 -> 4483.3845
$ws.Cells.Item(126, 12).Value2 = 37514817  # L126: 50016756 -> 37514817
$ws.Cells.Item(126, 13).Value2 = -2013.3845  # M126: -1606.5161 -> -2013.3845
$ws.Cells.Item(126, 14).Value2 = -37519757  # N126: -50021696 -> -37519757

Write-Output "applied 220 cell updates"